$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Append two new daily rows for 2025-12-24 (serial 46015):
#   row 48 -> 四方坪站 (shared string index 4)
#   row 49 -> 高岭站   (shared string index 5)
$ws.Cells.Item(48, 1).Value = 46015
$ws.Cells.Item(48, 2).Value = "四方坪站"
$ws.Cells.Item(48, 3).Value = 10754.29
$ws.Cells.Item(48, 4).Value = 9355.76
$ws.Cells.Item(48, 5).Value = 3508.97
$ws.Cells.Item(48, 6).Value = 439

$ws.Cells.Item(49, 1).Value = 46015
$ws.Cells.Item(49, 2).Value = "高岭站"
$ws.Cells.Item(49, 3).Value = 5147.93
$ws.Cells.Item(49, 4).Value = 4471.93
$ws.Cells.Item(49, 5).Value = 1375.42
$ws.Cells.Item(49, 6).Value = 182

# Mirror the view state Excel leaves behind after entering the new rows
$ws.Range("I49").Select() | Out-Null
